$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "cfr_respiratory_depression"
$ws.Range("B22").Value = 0.15

$ws.Range("A22").Select()
